# Update the build timestamp embedded in the "version" strings throughout the workbook.
# Old timestamp: January 30 2026 16.19.47 EST
# New timestamp: February 02 2026 12.49.33 EST

$oldVersion = "mines - January 30 (built on January 30 2026 16.19.47 EST)"
$newVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"

$wb = $excel.ActiveWorkbook

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---
$aboutSheet.Range("A2").Value = "Version: " + $newVersion

$aboutSheet.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Sanyuan Coal Mine, China, M2103, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet ---
# build_version column is S, rows 2 through 8 contain the old version string
for ($row = 2; $row -le 8; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    if ($cell.Value2 -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
